# ADD results from server
# Update the row-2 result values on each year sheet (2025, 2030, 2035, 2040, 2045, 2050)
# with the freshly computed figures coming back from the server run.

$wb = $excel.ActiveWorkbook

# NOTE: PowerShell hashtables coerce numeric-looking string keys (like "2025")
# to integers, which breaks Worksheets.Item(<name>) lookups. Use parallel
# arrays (sheet name + cell/value pairs) instead to keep the sheet names as
# real strings.

$sheetNames = @("2025", "2030", "2035", "2040", "2045", "2050")

$cellRefs = @("B2", "E2", "G2", "I2", "L2", "M2", "N2", "O2")

$valuesBySheet = @(
    @(1037.265132737054, 28926.05393052954, 8095.925712661834, 16171.06685703679, 48492.22142001599, 10595.37713982,    7015.544443014018, 6978.613354318873),
    @(4157.588990853394, 45991.90904307188, 8095.925712661834, 37079.12819938764, 54844.03303316472, 17449.04999683176, 8950.626290977361, 9689.183138434251),
    @(6368.910634126893, 57457.45307013817, 8095.925712661834, 52465.73681402855, 54844.03303316472, 21912.87293902603, 12955.24527970918, 12824.52739324476),
    @(6368.910634126893, 57457.45307013817, 8095.925712661834, 52465.73681402855, 54844.03303316472, 21912.87293902603, 13072.4393579009,  12824.52739324476),
    @(6368.910634126893, 57457.45307013817, 8095.925712661834, 52465.73681402855, 54844.03303316472, 21912.87293902603, 13521.65671205384, 14901.48629768362),
    @(6368.910634126893, 57457.45307013817, 8095.925712661834, 52465.73681402855, 54844.03303316472, 21912.87293902603, 13521.65671205384, 14901.48629768362)
)

# Column A2 only changes on the first two sheets (2025 -> 0, 2030 -> 0); the
# rest already hold their target value, but we set them explicitly so every
# sheet matches the server output exactly.
$aValues = @(0, 0, 2754.31755456332, 2754.31755456332, 5713.151062849596, 5713.151062849596)

for ($i = 0; $i -lt $sheetNames.Length; $i++) {
    $ws = $wb.Worksheets.Item($sheetNames[$i])

    $ws.Range("A2").Value = $aValues[$i]

    $rowValues = $valuesBySheet[$i]
    for ($j = 0; $j -lt $cellRefs.Length; $j++) {
        $ws.Range($cellRefs[$j]).Value = $rowValues[$j]
    }
}
